$d = $word.ActiveDocument

# --- Text content replacements -------------------------------------------------

# Title
$d.Content.Find.Execute("Biodiversity: Earth's Symphony of Life", $true, $false, $false, $false, $false, $true, 1, $false, "Exploring the Wonders of Mathematics: A Journey Through Numbers and Beyond", 2)

# Author name (collapses the 5 runs of "Dr" / "." / " Sarah J" / "." / " Carter" into one)
$d.Content.Find.Execute("Dr. Sarah J. Carter", $true, $false, $false, $false, $false, $true, 1, $false, "Evelyn Robinson", 2)

# Author e-mail handle/domain (keep the punctuation / "edu" runs separate)
$d.Content.Find.Execute("sarah", $true, $false, $false, $false, $false, $true, 1, $false, "evelyn", 2)
$d.Content.Find.Execute("carter@greenville", $true, $false, $false, $false, $false, $true, 1, $false, "robinson@highschool", 2)

# Body paragraph - 3 sentence-groups separated by line breaks
$d.Content.Find.Execute("Our planet, Earth, hosts an intricate tapestry of life forms, from microscopic organisms to majestic whales, weaving together a symphony of existence known as biodiversity. Each species, with its unique traits and interactions, contributes to the intricate web of life that sustains our planet. This spectacular array of living organisms is the result of billions of years of evolution, a process driven by natural selection. Biodiversity ensures the vitality and resilience of ecosystems, providing essential services like food production, water purification, and air filtration.", $true, $false, $false, $false, $false, $true, 1, $false, "In the vast realm of human knowledge and intellectual pursuit, Mathematics stands as a beacon of precision and elegance. For centuries, civilizations across the globe have harnessed the power of numbers and abstract concepts to make sense of the world around us. Whether it's calculating the movement of celestial bodies, unlocking the secrets of nature's patterns, or building intricate structures, Mathematics has been an indispensable tool, shaping our societies and propelling us into an era of scientific advancements.", 2)

$d.Content.Find.Execute("Furthermore, biodiversity enriches our lives in countless ways, inspiring art, music, and cultural traditions across human societies. It holds immense promise for scientific discovery, offering insights into the workings of the natural world and potential solutions to global challenges such as climate change and food security. Understanding and preserving biodiversity is crucial for ensuring a sustainable future for both human societies and the natural world.", $true, $false, $false, $false, $false, $true, 1, $false, "Like a grand tapestry woven with symbols and formulas, Mathematics invites us to explore a universe of possibilities. It is a language of logic and reason, enabling us to decipher complex phenomena and unravel hidden truths. From the earliest civilizations' counting systems to the sophisticated equations of the modern world, Mathematics has evolved as an essential element of human progress.", 2)

$d.Content.Find.Execute("In this essay, we will delve into the diverse realms of biodiversity, exploring how different species interact within ecosystems and uncovering the mechanisms that drive evolutionary change. We will also examine the threats facing biodiversity, such as habitat loss, pollution, and climate change, and discuss strategies for conserving this precious heritage for generations to come.", $true, $false, $false, $false, $false, $true, 1, $false, "At the heart of Mathematics lies the beauty of patterns and relationships. Whether it's the golden ratio found in art and architecture or the fractal patterns seen in nature, Mathematics provides a framework for understanding the interconnectedness of all things. It encourages us to think critically, analyze information, and recognize the underlying order within chaos.", 2)

# Summary heading paragraph stays "Summary" (unchanged)

# Summary paragraph - replace whole text (also drops the final two sentences)
$d.Content.Find.Execute("Biodiversity, the intricate tapestry of life on Earth, encompasses the vast array of species, their genetic variation, and the ecosystems they inhabit. It underpins the functioning of ecosystems, providing essential services for human well-being. Biodiversity enriches our lives culturally and holds promise for scientific discovery. However, human activities pose significant threats to biodiversity, leading to habitat loss, pollution, and climate change. Concerted efforts are needed to conserve biodiversity, including habitat protection, sustainable resource management, and education. Preserving biodiversity is vital for securing a sustainable future for both humanity and the natural world.", $true, $false, $false, $false, $false, $true, 1, $false, "In this essay, we embarked on an enthralling journey through the realm of Mathematics, exploring its precision, elegance, and its instrumental role in shaping human societies. We discovered Mathematics as a universal language, capable of describing the cosmos, unraveling nature's secrets, and constructing intricate structures. We marveled at the patterns and relationships revealed through mathematical inquiry, recognizing its ability to bring order to chaos and foster critical thinking skills. As we continue our exploration of Mathematics, we unlock new gateways to knowledge and deepen our understanding of the universe we inhabit.", 2)

# --- New trailing empty paragraph -----------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

# --- Font typo fix: TimesNewToman -> Times New Roman (applied last, per paragraph
#     so the paragraph-mark rPr is left untouched) --------------------------------
$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    $p = $d.Paragraphs($i)
    $s = $p.Range.Start
    $e = $p.Range.End
    if ($e -gt $s) {
        $r = $d.Range($s, $e - 1)
        $r.Font.Name = "Times New Roman"
    }
}
